$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54-170 down to 55-171
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly price record
$ws.Range("A54").Value = 11
$ws.Range("B54").Value = "Vega Monumental Concepción"
$ws.Range("C54").Value = "Bíobío"
$ws.Range("D54").Value = 44498
$ws.Range("D54").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E54").Value = 8
$ws.Range("F54").Value = 100112023
$ws.Range("G54").Value = "Brócoli"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 3000
$ws.Range("K54").Value = 650
$ws.Range("L54").Value = 700
$ws.Range("M54").Value = 675
$ws.Range("N54").Value = "$/unidad"
$ws.Range("O54").Value = "Región Metropolitana"
$ws.Range("P54").Value = 675
$ws.Range("Q54").Value = 1
$ws.Range("R54").Value = "Hortaliza"
